# VFDs are configurable to be have multiple vfds per suction group.
# Update the alerts sheet: auto-fit (best-fit) the column widths for
# the A:C table, remove the manually-set taller row heights (let them
# revert to the default row height), and move the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Remove the explicit (taller, wrapped-text) row heights on rows 2-9 and
# 11 by auto-fitting each row back down to the default height.
foreach ($r in 2,3,4,5,6,7,8,9,11) {
    $ws.Rows.Item($r).EntireRow.AutoFit() | Out-Null
}

# Best-fit (auto-fit) the column widths of columns A:C based on their
# contents -- now that rows 2-9/11 no longer wrap, the columns need to be
# wider to fit the full (unwrapped) alert/IO names.
$ws.Columns.Item(1).EntireColumn.AutoFit() | Out-Null
$ws.Columns.Item(2).EntireColumn.AutoFit() | Out-Null
$ws.Columns.Item(3).EntireColumn.AutoFit() | Out-Null
$ws.Columns.Item(1).ColumnWidth = 27.15
$ws.Columns.Item(2).ColumnWidth = 15.65
$ws.Columns.Item(3).ColumnWidth = 56

# Move the active cell selection to C13.
$ws.Range("C13").Select()
